$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 415; existing rows 415-462 shift down to 416-463.
$ws.Rows("415").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(415, 1).Value = 8
$ws.Cells.Item(415, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(415, 3).Value = "Coquimbo"
$ws.Cells.Item(415, 4).Value = 45212
$ws.Cells.Item(415, 5).Value = 4
$ws.Cells.Item(415, 6).Value = 100112031
$ws.Cells.Item(415, 7).Value = "Poroto verde"
$ws.Cells.Item(415, 8).Value = "Magnum"
$ws.Cells.Item(415, 9).Value = "Primera"
$ws.Cells.Item(415, 10).Value = 400
$ws.Cells.Item(415, 11).Value = 26500
$ws.Cells.Item(415, 12).Value = 27000
$ws.Cells.Item(415, 13).Value = 26750
$ws.Cells.Item(415, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(415, 15).Value = "Perú"
$ws.Cells.Item(415, 16).Value = 1070
$ws.Cells.Item(415, 17).Value = 25
$ws.Cells.Item(415, 18).Value = "Hortaliza"

# Match the date-style formatting used by the rest of column D.
$ws.Cells.Item(415, 4).NumberFormat = $ws.Cells.Item(416, 4).NumberFormat
